# Removal of unused samples.
#
# The following strain/sample rows are no longer used and must be removed
# entirely from the worksheet (row + the now-orphaned shared string):
#   b5145, b5158, b5161, b5164, b5189

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$samplesToRemove = @("b5145", "b5158", "b5161", "b5164", "b5189")

# Find the last used row in column A (strain names live there).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Collect the row numbers whose "strain" (column A) matches a sample to remove.
$rowsToDelete = @()
for ($r = 1; $r -le $lastRow; $r++) {
    $cellValue = $ws.Cells.Item($r, 1).Value2
    if ($samplesToRemove -contains $cellValue) {
        $rowsToDelete += $r
    }
}

# Delete from the bottom up so earlier row numbers stay valid as we go.
$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
